# Workbook: Hortaliza, Vega Monumental Concepción - Cilantro
# The commit adds a new week's worth of price observations (2 rows: "Primera"
# and "Segunda" quality) for the Cilantro series, inserted at row 89 so the
# whole table keeps growing downward (dimension goes from A1:R276 to A1:R278).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 89-90; everything from the old row 89 onward
# shifts down by two rows, preserving per-cell styles (e.g. the date style on
# column D).
$ws.Range("A89:A90").EntireRow.Insert()

# Constant columns shared by every record in this subset.
$mercadoId = 11
$mercado = "Vega Monumental Concepción"
$region = "Bíobío"
$codreg = 8
$categoriaId = 100112040
$categoria = "Cilantro"
$variedad = "Sin especificar"
$clasificacion = "Hortaliza"

# New row 89: "Primera" quality observation for 2023-02-10
$ws.Range("A89").Value = $mercadoId
$ws.Range("B89").Value = $mercado
$ws.Range("C89").Value = $region
$ws.Range("D89").Value = "2023-02-10"
$ws.Range("E89").Value = $codreg
$ws.Range("F89").Value = $categoriaId
$ws.Range("G89").Value = $categoria
$ws.Range("H89").Value = $variedad
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 700
$ws.Range("L89").Value = 800
$ws.Range("M89").Value = 750
$ws.Range("N89").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O89").Value = "Región de Ñuble"
$ws.Range("P89").Value = 750
$ws.Range("Q89").Value = 1
$ws.Range("R89").Value = $clasificacion

# New row 90: "Segunda" quality observation for the same date
$ws.Range("A90").Value = $mercadoId
$ws.Range("B90").Value = $mercado
$ws.Range("C90").Value = $region
$ws.Range("D90").Value = "2023-02-10"
$ws.Range("E90").Value = $codreg
$ws.Range("F90").Value = $categoriaId
$ws.Range("G90").Value = $categoria
$ws.Range("H90").Value = $variedad
$ws.Range("I90").Value = "Segunda"
$ws.Range("J90").Value = 100
$ws.Range("K90").Value = 600
$ws.Range("L90").Value = 600
$ws.Range("M90").Value = 600
$ws.Range("N90").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O90").Value = "Región de Ñuble"
$ws.Range("P90").Value = 600
$ws.Range("Q90").Value = 1
$ws.Range("R90").Value = $clasificacion
